$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction17")

$values = @(0,1,2,3,4,5,9,11,12,14,15,16,17,18,19,20,21,22,24,26,27,28,29,30,33,34,35,36,37,38,39,40,41,42,43,44,45,46,49,50,52,53,55,57,58,59,60,61,64,65,66,67,68,69,70,71,72,73,74,75)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $values[$i]
}
